$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; rows 51-94 shift down to 52-95
$ws.Rows("51:51").Insert()

# Populate the new row 51 with the data for the new record
$ws.Range("A51").Value = 11
$ws.Range("B51").Value = "Vega Monumental Concepción"
$ws.Range("C51").Value = "Bíobío"
$ws.Range("D51").Value2 = 44658
$ws.Range("E51").Value = 8
$ws.Range("F51").Value = 100112021
$ws.Range("G51").Value = "Ají"
$ws.Range("H51").Value = "Americana (o)"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 25
$ws.Range("K51").Value = 20000
$ws.Range("L51").Value = 22000
$ws.Range("M51").Value = 21200
$ws.Range("N51").Value = "$/caja 25 kilos"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 848
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
